$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MarchRaw")

# Header row: B1/C1/D1 use the "this month" shared-string headers matching
# the other *Raw sheets (Items owned by this library / other libraries / Total).
$janWs = $wb.Worksheets.Item("January")
$headerB1 = $janWs.Cells.Item(1, 2).Value2
$headerC1 = $janWs.Cells.Item(1, 3).Value2
$headerD1 = $janWs.Cells.Item(1, 4).Value2

$ws.Range("A1").Value = "Library"
$ws.Range("B1").Value = $headerB1
$ws.Range("C1").Value = $headerC1
$ws.Range("D1").Value = $headerD1

$rows = @(
    @(1, 'Library', $null, $null, $null),
    @(2, 'Atchison Public Library', 4465, 1358, 5823),
    @(3, 'Baldwin City Public Library', 2634, 619, 3253),
    @(4, 'Basehor Community Library', 8241, 1290, 9531),
    @(5, 'Bern Community Library', 148, 26, 174),
    @(6, 'Bonner Springs City Library', 5064, 1206, 6270),
    @(7, 'Burlingame Community Library', 509, 291, 800),
    @(8, 'Carbondale City Library', 624, 109, 733),
    @(9, 'Centralia Community Library', 252, 38, 290),
    @(10, 'Corning City Library', 91, $null, 91),
    @(11, 'Digital Content', $null, $null, $null),
    @(12, 'Doniphan County Library - Elwood', 48, 21, 69),
    @(13, 'Doniphan County Library - Highland', 325, 166, 491),
    @(14, 'Doniphan County Library - Troy', 504, 104, 608),
    @(15, 'Doniphan County Library - Wathena', 408, 44, 452),
    @(16, 'Effingham Community Library', 304, 70, 374),
    @(17, 'Eudora Community Library', 1748, 610, 2358),
    @(18, 'Everest, Barnes Reading Room', 141, 64, 205),
    @(19, 'Hiawatha, Morrill Public Library', 1842, 555, 2397),
    @(20, 'Highland Community College', 51, 2, 53),
    @(21, 'Holton, Beck-Bookman Library', 1811, 464, 2275),
    @(22, 'Horton Public Library', 116, 27, 143),
    @(23, 'Lansing Community Library', 1710, 636, 2346),
    @(24, 'Leavenworth Public Library', 9102, 2112, 11214),
    @(25, 'Linwood Community Library', 599, 168, 767),
    @(26, 'Louisburg Library', $null, $null, $null),
    @(27, 'Lyndon Carnegie Library', 354, 264, 618),
    @(28, 'McLouth Public Library', 185, 73, 258),
    @(29, 'Meriden-Ozawkie Public Library', 1298, 662, 1960),
    @(30, 'Northeast Kansas Library System', 20, 47, 67),
    @(31, 'Nortonville Public Library', 321, 52, 373),
    @(32, 'Osage City Library', 1342, 412, 1754),
    @(33, 'Osawatomie Public Library', 857, 299, 1156),
    @(34, 'Oskaloosa Public Library', 481, 182, 663),
    @(35, 'Ottawa Library', 6384, 865, 7249),
    @(36, 'Overbrook Public Library', 818, 220, 1038),
    @(37, 'Paola Free Library', 3069, 487, 3556),
    @(38, 'Perry-Lecompton Community Library', 61, 27, 88),
    @(39, 'Pomona Community Library', 56, 109, 165),
    @(40, 'Prairie Hills Schools - Axtell Public School', 335, 23, 358),
    @(41, 'Prairie Hills Schools - Sabetha Elementary School', 1571, 62, 1633),
    @(42, 'Prairie Hills Schools - Sabetha High School', 43, 3, 46),
    @(43, 'Prairie Hills Schools - Sabetha Middle School', 99, 9, 108),
    @(44, 'Prairie Hills Schools - Wetmore Academic Center (Permanently closed)', $null, $null, $null),
    @(45, 'Richmond Public Library', 346, 85, 431),
    @(46, 'Rossville Community Library', 1329, 444, 1773),
    @(47, 'Sabetha, Mary Cotton Library', 3061, 1079, 4140),
    @(48, 'Seneca Free Library', 1618, 299, 1917),
    @(49, 'Silver Lake Library', 958, 585, 1543),
    @(50, 'Tonganoxie Public Library', 3231, 924, 4155),
    @(51, 'Valley Falls, Delaware Township Library', 414, 271, 685),
    @(52, 'Wellsville City Library', 1142, 378, 1520),
    @(53, 'Wetmore Public Library', 114, 127, 241),
    @(54, 'Williamsburg Community Library', 245, 14, 259),
    @(55, 'Winchester Public Library', 328, 268, 596)
)

foreach ($row in $rows) {
    $r = $row[0]
    if ($r -eq 1) { continue }
    $name = $row[1]
    $b = $row[2]
    $c = $row[3]
    $d = $row[4]

    $ws.Cells.Item($r, 1).Value = $name
    if ($null -ne $b) { $ws.Cells.Item($r, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($r, 3).Value = $c }
    if ($null -ne $d) { $ws.Cells.Item($r, 4).Value = $d }
}

$excel.CalculateFullRebuild()
